$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The commit "Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta"
# removes the worker row for EDGAR RAFAEL BELTRAN MENDOZA (7931185 / period 2507,
# values 85410 / 2135250) that used to sit at row 16, leaving the remaining worker
# (SERGIO MANUEL DE AVILA HIDALGO, row 17) to shift up into row 16. It also
# refreshes the summary figures above the table.

$ws.Rows.Item(16).Delete()

# Update "Valor Mora" summary total
$ws.Range("E11").Value = 984

# Update worker / period counts (Cant. Trabajadores / Cant. Periodos)
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 1
